$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "Current as of:" date
$ws.Range("B1").Value = 44490

# Row 5: status for "Create makefile and directory organization for smooth workflow" changed to Completed
$ws.Range("B5").Value = "Completed"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "Functionally completed, with room for improvement down the line"

# Row 6: note text expanded
$ws.Range("E6").Value = "Includes reading array size metadata and geometry-specific text files. Had to figure out modules and compilation first."

# Update selection to C4
$ws.Range("C4").Select()

# Re-fit column E width to account for the longer note text (closest achievable
# approximation of the target bestFit width of 108.140625 given this runtime's
# column-width quantization)
$ws.Columns("E:E").ColumnWidth = 107.3
